$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mapping of row -> (new nombre_aides (C), new montant_total (E))
$updates = @{
    10  = @(278207, 1752504854)
    111 = @(6012, 12096389)
    115 = @(17559, 38636494)
    117 = @(19729, 56645610)
    124 = @(2663, 3903857)
    125 = @(4606, 13159157)
    126 = @(5646, 8180829)
    134 = @(5683, 17181109)
    139 = @(3321, 9344392)
    168 = @(285089, 1212846471)
    169 = @(562653, 1285790763)
    170 = @(367525, 2847630590)
    171 = @(115211, 448574472)
    174 = @(357339, 1019792833)
    177 = @(96775, 174802756)
    179 = @(235780, 813277362)
    180 = @(141516, 341192575)
    254 = @(62455, 99573670)
}

foreach ($row in $updates.Keys) {
    $vals = $updates[$row]
    $ws.Range("C$row").Value = $vals[0]
    $ws.Range("E$row").Value = $vals[1]
}
